{"js": "// Replace the date line and each \"a\u00f7b=\" expression in the table, in\n// document order, as described by the diff. A couple of the source\n// expressions (e.g. \"65\u00f77=\") repeat more than once in the document but\n// map to different results, so each occurrence is addressed by its\n// position among same-text matches (search results come back in\n// document order), not with a single global find/replace.\nconst replacements = [\n  [\"2023-11-02 Thursday\", \"2023-11-03 Friday\"],\n  [\"93\u00f73=\", \"24\u00f77=\"],\n  [\"65\u00f77=\", \"70\u00f75=\"],\n  [\"15\u00f72=\", \"47\u00f75=\"],\n  [\"49\u00f76=\", \"98\u00f75=\"],\n  [\"93\u00f76=\", \"44\u00f78=\"],\n  [\"49\u00f74=\", \"96\u00f77=\"],\n  [\"65\u00f77=\", \"15\u00f73=\"],\n  [\"92\u00f77=\", \"26\u00f75=\"],\n  [\"79\u00f72=\", \"41\u00f77=\"],\n  [\"27\u00f79=\", \"48\u00f76=\"],\n  [\"65\u00f76=\", \"16\u00f78=\"],\n  [\"46\u00f76=\", \"80\u00f73=\"],\n  [\"43\u00f79=\", \"14\u00f77=\"],\n  [\"35\u00f73=\", \"17\u00f72=\"],\n  [\"30\u00f76=\", \"28\u00f72=\"],\n  [\"29\u00f78=\", \"53\u00f79=\"],\n  [\"92\u00f78=\", \"48\u00f77=\"],\n  [\"57\u00f77=\", \"54\u00f77=\"],\n  [\"94\u00f79=\", \"15\u00f74=\"],\n  [\"58\u00f78=\", \"58\u00f75=\"],\n  [\"57\u00f75=\", \"76\u00f78=\"],\n  [\"69\u00f78=\", \"35\u00f79=\"],\n  [\"54\u00f74=\", \"64\u00f75=\"],\n  [\"96\u00f77=\", \"73\u00f78=\"],\n  [\"72\u00f73=\", \"87\u00f72=\"],\n];\n\n// Search once per distinct \"old\" text and cache the (ordered) ranges.\nconst searchCache = new Map();\nfor (const [oldText] of replacements) {\n  if (!searchCache.has(oldText)) {\n    const results = context.document.body.search(oldText, { matchCase: true });\n    results.load(\"items\");\n    searchCache.set(oldText, results);\n  }\n}\nawait context.sync();\n\n// Walk the replacements in document order, consuming one match per\n// duplicate text each time it is encountered.\nconst consumed = new Map();\nfor (const [oldText, newText] of replacements) {\n  const results = searchCache.get(oldText);\n  const idx = consumed.get(oldText) || 0;\n  const range = results.items[idx];\n  range.insertText(newText, \"Replace\");\n  consumed.set(oldText, idx + 1);\n}\nawait context.sync();\n", "ps1": "# Replace the date line and each \"a\u00f7b=\" expression in the table, in\n# document order, as described by the diff. A couple of the source\n# expressions (e.g. \"65\u00f77=\", and the \"96\u00f77=\" that one of the edits\n# itself introduces) repeat more than once across the whole edit, so a\n# single document-wide Find/Replace-All per pair is not safe - it would\n# also rewrite occurrences that must keep their original (different)\n# replacement, and a freshly-inserted replacement could accidentally be\n# re-matched by a later step searching for its old text.\n#\n# Instead we walk forward through the document exactly once: each\n# iteration searches only the remaining range (from the end of the\n# previous replacement to the end of the document), replaces just the\n# next single match (wdReplaceOne), and then advances the cursor past\n# what was just written before moving on to the next pair.\n$d = $word.ActiveDocument\n$cursor = 0\n\n$pairs = @(\n    @(\"2023-11-02 Thursday\", \"2023-11-03 Friday\"),\n    @(\"93\u00f73=\", \"24\u00f77=\"),\n    @(\"65\u00f77=\", \"70\u00f75=\"),\n    @(\"15\u00f72=\", \"47\u00f75=\"),\n    @(\"49\u00f76=\", \"98\u00f75=\"),\n    @(\"93\u00f76=\", \"44\u00f78=\"),\n    @(\"49\u00f74=\", \"96\u00f77=\"),\n    @(\"65\u00f77=\", \"15\u00f73=\"),\n    @(\"92\u00f77=\", \"26\u00f75=\"),\n    @(\"79\u00f72=\", \"41\u00f77=\"),\n    @(\"27\u00f79=\", \"48\u00f76=\"),\n    @(\"65\u00f76=\", \"16\u00f78=\"),\n    @(\"46\u00f76=\", \"80\u00f73=\"),\n    @(\"43\u00f79=\", \"14\u00f77=\"),\n    @(\"35\u00f73=\", \"17\u00f72=\"),\n    @(\"30\u00f76=\", \"28\u00f72=\"),\n    @(\"29\u00f78=\", \"53\u00f79=\"),\n    @(\"92\u00f78=\", \"48\u00f77=\"),\n    @(\"57\u00f77=\", \"54\u00f77=\"),\n    @(\"94\u00f79=\", \"15\u00f74=\"),\n    @(\"58\u00f78=\", \"58\u00f75=\"),\n    @(\"57\u00f75=\", \"76\u00f78=\"),\n    @(\"69\u00f78=\", \"35\u00f79=\"),\n    @(\"54\u00f74=\", \"64\u00f75=\"),\n    @(\"96\u00f77=\", \"73\u00f78=\"),\n    @(\"72\u00f73=\", \"87\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $r = $d.Range($cursor, $d.Content.End)\n    $find = $r.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 1) | Out-Null\n    $cursor = $r.End\n}\n"}
